$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "locus_name" column (B) for all data rows to the new constant value.
#    This also causes the now-unused shared string "ichthama.2-seasnakes" to be
#    dropped and the new one appended to the shared-string table automatically.
$ws.Range("B2:B7").Value = "ichthama.2-serpentes-UR"

# 2. Shade the whole data block (A2:AD7) with a white fill.
$ws.Range("A2:AD7").Interior.Color = 16777215

# 3. Resize / add column widths to match the new layout.
$ws.Columns("B").ColumnWidth = 26.83
$ws.Columns("E").ColumnWidth = 23.83
$ws.Columns("F").ColumnWidth = 17.83
$ws.Columns("AD").ColumnWidth = 17

# 4. Update the view: drop the old scrolled position / selection and select D14.
$ws.Range("D14").Select()
